$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProjectApplication")

# Row 2 updates:
# - Project ID (B2): 1 -> 2
# - Applicant NRIC (C2): "T2109876H" -> "S1234567A"
# - Date (F2): 45767.66396525463 -> 45768.54040121528
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = "S1234567A"
$ws.Range("F2").Value = 45768.54040121528
